$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1653
$ws1.Range("F3").Value = 9067
$ws1.Range("F4").Value = 109
$ws1.Range("F5").Value = 506
$ws1.Range("F7").Value = 1140
$ws1.Range("F8").Value = 194
$ws1.Range("F9").Value = 53
$ws1.Range("F10").Value = 87
$ws1.Range("F11").Value = 5830
$ws1.Range("F17").Value = 161
$ws1.Range("F19").Value = 20
$ws1.Range("F24").Value = 2711

# --- Sheet "演出" (performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19
$ws2.Range("F3").Value = 40

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1653
$ws4.Range("F3").Value = 9067
$ws4.Range("F4").Value = 109
$ws4.Range("F5").Value = 19
$ws4.Range("F6").Value = 506
$ws4.Range("F8").Value = 1140
$ws4.Range("F9").Value = 194
$ws4.Range("F10").Value = 53
$ws4.Range("F11").Value = 87
$ws4.Range("F12").Value = 5830
$ws4.Range("F18").Value = 161
$ws4.Range("F20").Value = 20
$ws4.Range("F25").Value = 2711
$ws4.Range("F26").Value = 40
